{"js": "// Load all paragraphs in the document body.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\n// Helper: find a paragraph whose (concatenated) text starts with the given\n// snippet and replace a specific sub-string inside it with new text,\n// collapsing whichever runs span that sub-string into a single run.\n// This mirrors what happens when a user selects that exact span of text\n// in Word and retypes it.\nasync function replaceSpan(paragraph, oldSpan, newSpan) {\n  const results = paragraph.search(oldSpan, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Span not found: \" + oldSpan);\n  }\n  results.items[0].insertText(newSpan, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1) \"Con estos $4 no conviene comprar nada...\" -> \"Debo analizar la tabla\n//    optima alternativa para ver que hago con estos $4.\"\nconst p1 = paragraphs.items.find((p) =>\n  p.text.indexOf(\"Con estos $4 no conviene comprar nada\") !== -1\n);\np1.insertText(\n  \"Debo analizar la tabla optima alternativa para ver que hago con estos $4.\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// 2) \"4/3 * 17 + 12 * (-5/3) - 4 = \" (many runs) -> single run, keeping the\n//    bold \"-4/3 <= 0\" that follows untouched.\nparagraphs.load(\"items,text\");\nawait context.sync();\nconst p2 = paragraphs.items.find((p) =>\n  p.text.indexOf(\"4/3 * 17 + 12 * (-5/3) - 4 = \") !== -1\n);\nawait replaceSpan(p2, \"4/3 * 17 + 12 * (-5/3) - 4 = \", \"4/3 * 17 + 12 * (-5/3) - 4 = \");\n\n// 3) \"-2/3 * 17 + 12 * 1/3 = \" (many runs) -> single run, bold \"-22/3 <= 0\"\n//    left untouched.\nparagraphs.load(\"items,text\");\nawait context.sync();\nconst p3 = paragraphs.items.find((p) =>\n  p.text.indexOf(\"-2/3 * 17 + 12 * 1/3 = \") !== -1\n);\nawait replaceSpan(p3, \"-2/3 * 17 + 12 * 1/3 = \", \"-2/3 * 17 + 12 * 1/3 = \");\n\n// 4) \"1/3* 17 + 12 * (-2/3) = \" (many runs) -> single run, bold \"-7/3 <= 0\"\n//    left untouched.\nparagraphs.load(\"items,text\");\nawait context.sync();\nconst p4 = paragraphs.items.find((p) =>\n  p.text.indexOf(\"1/3* 17 + 12 * (-2/3) = \") !== -1\n);\nawait replaceSpan(p4, \"1/3* 17 + 12 * (-2/3) = \", \"1/3* 17 + 12 * (-2/3) = \");\n\n// 5) \"4/3 * 12 + 12 * (-5/3) \u2013 14 = -18\" (many runs spanning the whole\n//    paragraph) -> single run.\nparagraphs.load(\"items,text\");\nawait context.sync();\nconst p5 = paragraphs.items.find((p) =>\n  p.text.indexOf(\"4/3 * 12 + 12 * (-5/3)\") !== -1 && p.text.indexOf(\"14 = -18\") !== -1\n);\np5.insertText(\"4/3 * 12 + 12 * (-5/3) \u2013 14 = -18\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 6) \"Si compro la m\u00e1quina b estar\u00eda ganando $16 por mes -> \" (several\n//    runs) -> single run; \"anualmente\" / \" $\" / \"192\" runs stay untouched.\nparagraphs.load(\"items,text\");\nawait context.sync();\nconst p6 = paragraphs.items.find((p) =>\n  p.text.indexOf(\"Si compro la m\u00e1quina b estar\u00eda ganando $16 por mes\") !== -1\n);\nawait replaceSpan(\n  p6,\n  \"Si compro la m\u00e1quina b estar\u00eda ganando $16 por mes -> \",\n  \"Si compro la m\u00e1quina b estar\u00eda ganando $16 por mes -> \"\n);\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-Span {\n    param(\n        [string]$OldText,\n        [string]$NewText\n    )\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $result = $find.Execute(\n        $OldText,   # FindText\n        $true,      # MatchCase\n        $false,     # MatchWholeWord\n        $false,     # MatchWildcards\n        $false,     # MatchSoundsLike\n        $false,     # MatchAllWordForms\n        $true,      # Forward\n        1,          # Wrap (wdFindContinue)\n        $false,     # Format\n        $NewText,   # ReplaceWith\n        2           # Replace (wdReplaceAll)\n    )\n    if (-not $result) {\n        throw \"Could not find span: $OldText\"\n    }\n}\n\n# 1) \"Con estos $4 no conviene comprar nada...\" -> \"Debo analizar la tabla\n#    optima alternativa para ver que hago con estos $4.\"\n$d.Paragraphs(21).Range.Text = \"Debo analizar la tabla optima alternativa para ver que hago con estos `$4.\"\n\n# 2) \"4/3 * 17 + 12 * (-5/3) - 4 = \" (several runs) -> single run; the bold\n#    \"-4/3 <= 0\" that follows is left untouched.\nReplace-Span -OldText \"4/3 * 17 + 12 * (-5/3) - 4 = \" -NewText \"4/3 * 17 + 12 * (-5/3) - 4 = \"\n\n# 3) \"-2/3 * 17 + 12 * 1/3 = \" (several runs) -> single run; the bold\n#    \"-22/3 <= 0\" that follows is left untouched.\nReplace-Span -OldText \"-2/3 * 17 + 12 * 1/3 = \" -NewText \"-2/3 * 17 + 12 * 1/3 = \"\n\n# 4) \"1/3* 17 + 12 * (-2/3) = \" (several runs) -> single run; the bold\n#    \"-7/3 <= 0\" that follows is left untouched.\nReplace-Span -OldText \"1/3* 17 + 12 * (-2/3) = \" -NewText \"1/3* 17 + 12 * (-2/3) = \"\n\n# 5) \"4/3 * 12 + 12 * (-5/3) [en dash] 14 = -18\" (several runs spanning the\n#    whole paragraph) -> single run.\n$d.Paragraphs(34).Range.Text = \"4/3 * 12 + 12 * (-5/3) `u{2013} 14 = -18\"\n\n# 6) \"Si compro la m\u00e1quina b estar\u00eda ganando $16 por mes -> \" (several\n#    runs) -> single run; \"anualmente\" / \" $\" / \"192\" runs stay untouched.\nReplace-Span -OldText \"Si compro la m\u00e1quina b estar\u00eda ganando `$16 por mes -> \" -NewText \"Si compro la m\u00e1quina b estar\u00eda ganando `$16 por mes -> \"\n"}
